$d = $word.ActiveDocument

# 1) Fix the Portuguese typo "sequêncial" -> "sequencial" in the
#    "Finalmente, a versão não sequêncial..." paragraph.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("sequêncial foi de fácil", $true, $false, $false, $false, $false, `
               $true, 1, $false, "sequencial foi de fácil", 2)

# 2) Move the hidden "_GoBack" bookmark (Word's "last edit location"
#    marker) from its old spot - at the end of the following paragraph,
#    after "...estabelecer uma relação entre eles." - to right after the
#    "seque" that was just retyped (i.e. exactly where the "ê" -> "e"
#    correction was made), splitting the run the way Word does when it
#    drops the bookmark mid-run.

# Locate "Finalmente, a vers" to get the first split point.
$findPrefix = $d.Content.Find
$findPrefix.ClearFormatting()
$findPrefix.Text = "Finalmente, a vers"
$findPrefix.Execute()
$prefixEnd = $findPrefix.Parent.End

# Locate "Finalmente, a versão não seque" to get the bookmark's final
# position (right after the corrected "e").
$findAnchor = $d.Content.Find
$findAnchor.ClearFormatting()
$findAnchor.Text = "Finalmente, a versão não seque"
$findAnchor.Execute()
$anchorEnd = $findAnchor.Parent.End

# First drop the bookmark at the old intermediate split point (this is
# how the run ends up split into three pieces), then relocate it (a
# bookmark name is unique, so re-adding it moves it) to its real final
# resting place.
$rTmp = $d.Range($prefixEnd, $prefixEnd)
$d.Bookmarks.Add("_GoBack", $rTmp)

$rFinal = $d.Range($anchorEnd, $anchorEnd)
$d.Bookmarks.Add("_GoBack", $rFinal)
